# Replace "OIE" with "WOAH" in the specific cells that reference the
# World Organisation for Animal Health (OIE -> WOAH rename), while leaving
# unrelated content (e.g. oie.int URLs) untouched.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet 1")
$sheet2 = $wb.Worksheets.Item("References")

$sheet1.Range("E5").Value  = "Based on official disease reports to the WOAH"
$sheet1.Range("E6").Value  = "ASF is a disease listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code and must be reported to the WOAH. The map to the right displays outbreak points reported to the WOAH early warning system since 2005."
$sheet1.Range("E7").Value  = "As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:"
$sheet1.Range("E14").Value = "Countries are coloured according to the available information regarding their stable disease situation (disease status legend). This information is provided by countries through the WOAH monitoring system, which is a different reporting channel.<br>Immediate notifications (points) and disease status (country/region colours) are reported to the WOAH in different spatial and temporal scales, and therefore are displayed in the map as layers which can be filtered independently."
$sheet1.Range("E17").Value = "For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}."
$sheet1.Range("E31").Value = "A summary of the disease in animal hosts is given in the {ref008:WOAH Technical disease card}."
$sheet1.Range("E52").Value = "Refer to the {ref008:WOAH Technical disease card} for a key summary of the virus characteristics. "
$sheet1.Range("E64").Value = "Refer to the {ref008:WOAH Technical disease card} for a key summary of the disease transmission and epidemiological parameters."
$nbsp = [char]0x00A0
$sheet1.Range("E76").Value = "WOAH-prescribed tests for agent identification are: Virus isolation, fluorescent antibody test (FAT), ELISA for antigen detection and PCR. For the detection of immune response: ELISA, indirect immunoperoxidase test (IPT), indirect fluorescent antibody test (IFAT) and  immunoblotting test (IBT) ({ref010:WOAH," + $nbsp + "Terrestrial Manual})."
$sheet1.Range("E126").Value = "Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data."

$sheet2.Range("C2").Value  = "WOAH-WAHIS (WOAH World Animal Health Information System)"
$sheet2.Range("C5").Value  = "WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$sheet2.Range("C8").Value  = "WOAH (World Organisation for Animal Health) Technical Disease Card: African swine fever. 2021."
$sheet2.Range("C9").Value  = "WOAH (World Organisation for Animal Health), 2021. African Swine fever. Chapter 15.1. WOAH Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$sheet2.Range("C10").Value = "WOAH (World Organisation for Animal Health), 2019. African Swine fever. Chapter 3.08.01. WOAH Terrestrial Manual 2019. WOAH, Paris, France"
